# Fixed a bug in winResultModifier
# The rows of result data got shuffled; restore the correct
# row -> data association described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2  = @(101, 9, 30, 15, 60, 15)
    3  = @(1001, 18, 30, 75, 60, 72)
    4  = @(501, 9, 52, 30, 75, 45)
    7  = @(601, 9, 60, 67, 60, 42)
    8  = @(1201, 2, 10, 10, 10, 10)
    9  = @(701, 3, 90, 45, 97, 15)
    10 = @(201, 9, 30, 15, 45, 30)
    11 = @(901, 16, 15, 45, 60, 60)
    12 = @(301, 6, 45, 30, 60, 45)
    13 = @(801, 3, 67, 65, 52, 45)
    14 = @(1202, 2, 10, 10, 10, 10)
    15 = @(1203, 3, 15, 15, 15, 15)
    17 = @(2, 0, 2, 2, 2, 2)
    18 = @(3, 0, 3, 3, 3, 3)
    19 = @(1101, 0, 15, 30, 30, 0)
    20 = @(802, 0, 4, 5, 4, 0)
    21 = @(1, 0, 2, 2, 2, 2)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
